$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert new columns so we have room for the new "apoio_*" and
#    "contribuicoes_*" statistics, shifting the existing columns to the right.
# ---------------------------------------------------------------------------

# Insert 3 blank columns right after K (apoio_medio) -> they will hold
# apoio_std, apoio_min, apoio_max. This pushes the old L (contribuicoes) and
# M (media_contribuicoes) columns to O and P respectively.
$ws.Range("L1:N1").EntireColumn.Insert()

# Insert 3 more blank columns right after the (now shifted) P column
# (media_contribuicoes) -> they will hold contribuicoes_std, contribuicoes_min
# and contribuicoes_max. This pushes menor_ano/maior_ano to T and U.
$ws.Range("Q1:S1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2) Rename / retitle headers to match the new naming scheme.
# ---------------------------------------------------------------------------

$ws.Range("G1").Value = "arrecadado_avg"
$ws.Range("H1").Value = "arrecadado_std"
$ws.Range("I1").Value = "arrecadado_min"
$ws.Range("J1").Value = "arrecadado_max"

$ws.Range("L1").Value = "apoio_std"
$ws.Range("M1").Value = "apoio_min"
$ws.Range("N1").Value = "apoio_max"

$ws.Range("P1").Value = "contribuicoes_med"
$ws.Range("Q1").Value = "contribuicoes_std"
$ws.Range("R1").Value = "contribuicoes_min"
$ws.Range("S1").Value = "contribuicoes_max"

# ---------------------------------------------------------------------------
# 3) Update / fill in the data values for rows 2-4.
# ---------------------------------------------------------------------------

# Row 2 - aon
$ws.Range("K2").Value = 91.85574933975617
$ws.Range("L2").Value = 49.08980856017526
$ws.Range("M2").Value = 13.93896149503088
$ws.Range("N2").Value = 792.0360759681182
$ws.Range("Q2").Value = 423.019225146675
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 6494

# Row 3 - flex
$ws.Range("K3").Value = 77.41063997458096
$ws.Range("L3").Value = 39.50983355883143
$ws.Range("M3").Value = 10.77163914429046
$ws.Range("N3").Value = 461.5197709071476
$ws.Range("Q3").Value = 327.6748910926806
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 7954

# Row 4 - sub
$ws.Range("K4").Value = 21.28348419490777
$ws.Range("L4").Value = 15.01968006252796
$ws.Range("M4").Value = 1.011042153300025
$ws.Range("N4").Value = 84.0771316599004
$ws.Range("Q4").Value = 31.86830254134198
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 208

Write-Host "Edit applied"
